$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(-19.4016378655953,1.962484496384981,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953),
    @(-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,1.264795070097932,-19.4016378655953,-19.4016378655953),
    @(-19.4016378655953,1.967510549383099,1.685175422961444,-19.4016378655953,3.399296979542475,-19.4016378655953,1.385745832346961,-19.4016378655953,0.8148723704245016,-19.4016378655953),
    @(-19.4016378655953,1.705283480545667,-19.4016378655953,-19.4016378655953,-19.4016378655953,2.913112697899164,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953),
    @(-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953),
    @(2.43915119688107,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953),
    @(-19.4016378655953,-19.4016378655953,-19.4016378655953,1.793001777595635,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953),
    @(3.865597028967233,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953),
    @(-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,1.698432185193016,-19.4016378655953,4.321926115989044),
    @(-19.4016378655953,-19.4016378655953,-19.4016378655953,2.915980379477737,-19.4016378655953,2.938436334287047,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953),
    @(-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953),
    @(-19.4016378655953,-19.4016378655953,-19.4016378655953,2.517863084518607,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,1.576059211120866,-19.4016378655953),
    @(-19.4016378655953,-19.4016378655953,1.530320899965772,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953),
    @(-19.4016378655953,-19.4016378655953,1.75330444217075,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953),
    @(-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,1.939088617313267,-19.4016378655953),
    @(-19.4016378655953,2.137983397458781,1.826758685617999,-19.4016378655953,-19.4016378655953,-19.4016378655953,2.029539873757906,2.147088578082583,2.479690100763979,-19.4016378655953),
    @(-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,-19.4016378655953,2.068965903153444,2.020141912722466,2.54766772118136,-19.4016378655953),
    @(-19.4016378655953,-19.4016378655953,2.049402660810946,-19.4016378655953,-19.4016378655953,-19.4016378655953,1.595420986851112,1.882667198410783,-19.4016378655953,-19.4016378655953),
    @(-19.4016378655953,1.042431688904978,1.504695673864715,-19.4016378655953,3.240169809163277,-19.4016378655953,1.758214382185231,1.12329499307709,-19.4016378655953,-19.4016378655953),
    @(-19.4016378655953,1.304299044329603,-19.4016378655953,1.704823478429642,-19.4016378655953,2.263563143107639,1.434928011117364,-19.4016378655953,-19.4016378655953,-19.4016378655953)
)

for ($r = 0; $r -lt $newValues.Length; $r++) {
    $rowVals = $newValues[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
